# Updates recalculated market-price-derived figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) for several Leve rows across the ALC, ARM, BSM,
# CRP, GSM, LTW and WVR sheets, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4457574.5
$ws.Range("J17").Value = 4457574.5
$ws.Range("L17").Value = 13372723.5
$ws.Range("N17").Value = -13373059.5
$ws.Range("H64").Value = 5809.091
$ws.Range("I64").Value = 6277.778
$ws.Range("J64").Value = 3700
$ws.Range("K64").Value = 6277.778
$ws.Range("L64").Value = 3700
$ws.Range("M64").Value = -6029.778
$ws.Range("N64").Value = -4196
$ws.Range("H67").Value = 5809.091
$ws.Range("I67").Value = 6277.778
$ws.Range("J67").Value = 3700
$ws.Range("K67").Value = 6277.778
$ws.Range("L67").Value = 3700
$ws.Range("M67").Value = -5419.778
$ws.Range("N67").Value = -5416
$ws.Range("H76").Value = 4173072
$ws.Range("I76").Value = 11120838
$ws.Range("J76").Value = 4412.32
$ws.Range("K76").Value = 11120838
$ws.Range("L76").Value = 4412.32
$ws.Range("M76").Value = -11120523
$ws.Range("N76").Value = -5042.32
$ws.Range("H79").Value = 4173072
$ws.Range("I79").Value = 11120838
$ws.Range("J79").Value = 4412.32
$ws.Range("K79").Value = 11120838
$ws.Range("L79").Value = 4412.32
$ws.Range("M79").Value = -11119746
$ws.Range("N79").Value = -6596.32
$ws.Range("H112").Value = 13606521
$ws.Range("I112").Value = 325
$ws.Range("J112").Value = 15038752
$ws.Range("K112").Value = 975
$ws.Range("L112").Value = 45116256
$ws.Range("M112").Value = 133
$ws.Range("N112").Value = -45118472
$ws.Range("H116").Value = 9087.1
$ws.Range("I116").Value = 12814.546
$ws.Range("K116").Value = 12814.546
$ws.Range("M116").Value = -9372.546
$ws.Range("H132").Value = 2045.2258
$ws.Range("I132").Value = 1421.4117
$ws.Range("J132").Value = 4937.4546
$ws.Range("K132").Value = 4264.2351
$ws.Range("L132").Value = 14812.3638
$ws.Range("M132").Value = -1734.2351
$ws.Range("N132").Value = -19872.3638
$ws.Range("H137").Value = 2085
$ws.Range("I137").Value = 1591.3077
$ws.Range("J137").Value = 3251.9092
$ws.Range("K137").Value = 4773.9231
$ws.Range("L137").Value = 9755.7276
$ws.Range("M137").Value = -2223.9231
$ws.Range("N137").Value = -14855.7276
$ws.Range("H141").Value = 2903.9583
$ws.Range("I141").Value = 1846.2667
$ws.Range("J141").Value = 4666.778
$ws.Range("K141").Value = 5538.800099999999
$ws.Range("L141").Value = 14000.334
$ws.Range("M141").Value = -358.8000999999995
$ws.Range("N141").Value = -24360.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6708.1914
$ws.Range("I32").Value = 6005.8667
$ws.Range("J32").Value = 22510.5
$ws.Range("K32").Value = 6005.8667
$ws.Range("L32").Value = 22510.5
$ws.Range("M32").Value = -5718.8667
$ws.Range("N32").Value = -23084.5
$ws.Range("H44").Value = 17757.143
$ws.Range("I44").Value = 8266.666999999999
$ws.Range("J44").Value = 24875
$ws.Range("K44").Value = 8266.666999999999
$ws.Range("L44").Value = 24875
$ws.Range("M44").Value = -7778.666999999999
$ws.Range("N44").Value = -25851
$ws.Range("H55").Value = 19995.334
$ws.Range("J55").Value = 19995.334
$ws.Range("L55").Value = 19995.334
$ws.Range("N55").Value = -20625.334
$ws.Range("H63").Value = 333335330
$ws.Range("I63").Value = 333335330
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 333335330
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -333334644
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 333335330
$ws.Range("I66").Value = 333335330
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 1666676650
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -1666673218
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 2118.7917
$ws.Range("I74").Value = 1694.8667
$ws.Range("J74").Value = 2825.3333
$ws.Range("K74").Value = 1694.8667
$ws.Range("L74").Value = 2825.3333
$ws.Range("M74").Value = -820.8667
$ws.Range("N74").Value = -4573.3333
$ws.Range("H77").Value = 2118.7917
$ws.Range("I77").Value = 1694.8667
$ws.Range("J77").Value = 2825.3333
$ws.Range("K77").Value = 8474.333500000001
$ws.Range("L77").Value = 14126.6665
$ws.Range("M77").Value = -4106.333500000001
$ws.Range("N77").Value = -22862.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1990.6923
$ws.Range("I86").Value = 1922.1111
$ws.Range("K86").Value = 1922.1111
$ws.Range("M86").Value = -799.1111000000001
$ws.Range("H89").Value = 1990.6923
$ws.Range("I89").Value = 1922.1111
$ws.Range("K89").Value = 9610.5555
$ws.Range("M89").Value = -3994.5555
$ws.Range("H94").Value = 1229.7646
$ws.Range("I94").Value = 851.43475
$ws.Range("J94").Value = 2020.8182
$ws.Range("K94").Value = 851.43475
$ws.Range("L94").Value = 2020.8182
$ws.Range("M94").Value = -400.43475
$ws.Range("N94").Value = -2922.8182
$ws.Range("H105").Value = 2201.2666
$ws.Range("I105").Value = 2159.1538
$ws.Range("K105").Value = 2159.1538
$ws.Range("M105").Value = -412.1538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5771.5415
$ws.Range("I31").Value = 2758.6956
$ws.Range("J31").Value = 8543.360000000001
$ws.Range("K31").Value = 2758.6956
$ws.Range("L31").Value = 8543.360000000001
$ws.Range("M31").Value = -2463.6956
$ws.Range("N31").Value = -9133.360000000001
$ws.Range("H34").Value = 5771.5415
$ws.Range("I34").Value = 2758.6956
$ws.Range("J34").Value = 8543.360000000001
$ws.Range("K34").Value = 2758.6956
$ws.Range("L34").Value = 8543.360000000001
$ws.Range("M34").Value = -2556.6956
$ws.Range("N34").Value = -8947.360000000001
$ws.Range("H62").Value = 6875.25
$ws.Range("I62").Value = 6833.1665
$ws.Range("K62").Value = 6833.1665
$ws.Range("M62").Value = -6209.1665
$ws.Range("H65").Value = 6875.25
$ws.Range("I65").Value = 6833.1665
$ws.Range("K65").Value = 34165.8325
$ws.Range("M65").Value = -31045.8325

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5565.9443
$ws.Range("I70").Value = 5624.9165
$ws.Range("J70").Value = 5448
$ws.Range("K70").Value = 5624.9165
$ws.Range("L70").Value = 5448
$ws.Range("M70").Value = -5354.9165
$ws.Range("N70").Value = -5988
$ws.Range("H73").Value = 5565.9443
$ws.Range("I73").Value = 5624.9165
$ws.Range("J73").Value = 5448
$ws.Range("K73").Value = 5624.9165
$ws.Range("L73").Value = 5448
$ws.Range("M73").Value = -4688.9165
$ws.Range("N73").Value = -7320
$ws.Range("H80").Value = 6561.875
$ws.Range("I80").Value = 9463.214
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 9463.214
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -8465.214
$ws.Range("N80").Value = -4496
$ws.Range("H83").Value = 6561.875
$ws.Range("I83").Value = 9463.214
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 47316.07
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -42324.07
$ws.Range("N83").Value = -22484

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1944.6111
$ws.Range("J22").Value = 2029.5294
$ws.Range("L22").Value = 2029.5294
$ws.Range("N22").Value = -2619.5294
$ws.Range("H27").Value = 1944.6111
$ws.Range("J27").Value = 2029.5294
$ws.Range("L27").Value = 2029.5294
$ws.Range("N27").Value = -2243.5294
$ws.Range("H122").Value = 5435554.5
$ws.Range("I122").Value = 5960443
$ws.Range("K122").Value = 17881329
$ws.Range("M122").Value = -17878879
$ws.Range("H132").Value = 12828170
$ws.Range("I132").Value = 15881091
$ws.Range("J132").Value = 5900.8
$ws.Range("K132").Value = 47643273
$ws.Range("L132").Value = 17702.4
$ws.Range("M132").Value = -47640743
$ws.Range("N132").Value = -22762.4
$ws.Range("H136").Value = 8052.7046
$ws.Range("I136").Value = 5815.1665
$ws.Range("J136").Value = 12847.429
$ws.Range("K136").Value = 17445.4995
$ws.Range("L136").Value = 38542.287
$ws.Range("M136").Value = -14895.4995
$ws.Range("N136").Value = -43642.287

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2176.9138
$ws.Range("I136").Value = 2099.8918
$ws.Range("K136").Value = 6299.6754
$ws.Range("M136").Value = -3749.6754

Write-Output "Updated leve-profit figures on ALC, ARM, BSM, CRP, GSM, LTW, WVR"
